# Tutorial 6 solution update: switch the Date column from slash-separated
# to dash-separated formatting, and refresh the computed attendance
# counters that changed as a result of the updated source data.
#
# The date cells hold plain text (e.g. "28/07/2022"), not real Excel
# dates. Some of the new dash-separated strings (day <= 12) are
# ambiguous and Excel's COM layer would otherwise auto-convert them to
# date serials on assignment (e.g. "01-08-2022" -> 01 Aug 2022 date).
# Forcing the cell to Text before the write, then handing the style
# back to "Normal" afterwards, keeps the value a literal string without
# leaving any lingering number-format override on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date text: DD/MM/YYYY -> DD-MM-YYYY (rows 3..21) ---
Set-TextValue $ws.Range("A3")  "28-07-2022"
Set-TextValue $ws.Range("A4")  "01-08-2022"
Set-TextValue $ws.Range("A5")  "04-08-2022"
Set-TextValue $ws.Range("A6")  "08-08-2022"
Set-TextValue $ws.Range("A7")  "11-08-2022"
Set-TextValue $ws.Range("A8")  "15-08-2022"
Set-TextValue $ws.Range("A9")  "18-08-2022"
Set-TextValue $ws.Range("A10") "22-08-2022"
Set-TextValue $ws.Range("A11") "25-08-2022"
Set-TextValue $ws.Range("A12") "29-08-2022"
Set-TextValue $ws.Range("A13") "01-09-2022"
Set-TextValue $ws.Range("A14") "05-09-2022"
Set-TextValue $ws.Range("A15") "08-09-2022"
Set-TextValue $ws.Range("A16") "12-09-2022"
Set-TextValue $ws.Range("A17") "15-09-2022"
Set-TextValue $ws.Range("A18") "19-09-2022"
Set-TextValue $ws.Range("A19") "22-09-2022"
Set-TextValue $ws.Range("A20") "26-09-2022"
Set-TextValue $ws.Range("A21") "29-09-2022"

# --- Attendance counters that changed for row 3 (28-07-2022) ---
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# --- Attendance counters that changed for row 4 (01-08-2022) ---
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# --- Attendance counters that changed for row 5 (04-08-2022) ---
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# --- Attendance counters that changed for row 13 (01-09-2022) ---
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0
